# Adds the 2020-03-29 row (row 91) to both the "Confirmed" and "Deaths" sheets,
# mirroring the ECDC daily-totals layout already present in rows 2-90.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)

# Column A: date label "2020-03-29" — force text (avoid auto date-serial coercion)
# by entering it quote-prefixed, the way Excel treats a manually typed 'date, then
# clear the resulting quote-prefix style so the cell keeps the sheet's default format.
$ws.Cells.Item(91, 1).Value = "'2020-03-29"
$ws.Cells.Item(91, 1).Style = "Normal"

$row91 = @{
    2=106; 3=197; 4=409; 5=308; 6=4; 7=2; 8=7; 9=745; 10=424; 11=28;
    12=3809; 13=8291; 14=182; 15=11; 16=473; 17=48; 18=26; 19=94; 20=9134; 21=2;
    22=6; 23=22; 24=4; 25=81; 26=257; 27=3904; 28=2; 29=120; 30=331; 31=180;
    32=103; 33=99; 34=5386; 35=6; 37=8; 38=6; 39=5; 40=1909; 41=82342; 42=608;
    43=19; 44=295; 45=140; 46=657; 47=119; 48=8; 49=179; 50=2663; 51=58; 52=2201;
    53=15; 54=11; 55=719; 56=1835; 57=536; 58=19; 59=13; 60=6; 61=640; 62=9;
    63=16; 64=155; 65=5; 66=1218; 67=37575; 68=34; 69=7; 70=3; 71=85; 72=52547;
    73=141; 74=56; 75=1061; 76=10; 77=9; 78=55; 79=34; 80=39; 81=8; 82=2;
    83=8; 84=8; 85=6; 86=110; 87=408; 88=963; 89=979; 90=1155; 91=35408; 92=506;
    93=2415; 94=32; 95=3619; 96=92472; 97=32; 98=1693; 99=61; 100=235; 101=229; 102=38;
    103=86; 104=235; 105=84; 106=6; 107=305; 108=412; 109=3; 110=1; 111=61; 112=394;
    113=1831; 114=28; 115=2320; 116=16; 117=9; 118=139; 119=5; 120=102; 121=848; 122=231;
    123=43; 124=12; 125=84; 126=5; 127=358; 128=8; 129=8; 130=8; 131=5; 132=9762;
    133=15; 134=476; 135=3; 136=10; 137=97; 138=241; 139=3845; 140=152; 141=1408; 142=97;
    143=901; 144=1; 145=59; 146=671; 147=1075; 148=1638; 149=5170; 150=100; 151=590; 152=1452;
    153=1264; 154=60; 155=2; 156=4; 157=1; 158=224; 159=1203; 160=130; 161=659; 162=7;
    163=803; 164=3; 165=295; 166=691; 167=3; 168=1187; 169=9583; 170=72248; 171=115; 172=5;
    173=8; 174=3447; 175=13152; 176=5; 177=283; 178=1245; 179=1; 180=28; 181=74; 182=227;
    183=7402; 184=5; 185=30; 186=311; 187=468; 188=17089; 189=13; 190=22; 191=124665; 192=304;
    193=133; 194=119; 195=223; 196=28; 197=7;
}
foreach ($col in $row91.Keys) {
    $ws.Cells.Item(91, [int]$col).Value = $row91[$col]
}

$ws = $wb.Worksheets.Item(2)

# Column A: date label "2020-03-29" — force text (avoid auto date-serial coercion)
# by entering it quote-prefixed, the way Excel treats a manually typed 'date, then
# clear the resulting quote-prefix style so the cell keeps the sheet's default format.
$ws.Cells.Item(91, 1).Value = "'2020-03-29"
$ws.Cells.Item(91, 1).Style = "Normal"

$row91 = @{
    2=3; 3=10; 4=26; 5=4; 6=0; 7=0; 8=0; 9=19; 10=3; 11=0;
    12=14; 13=68; 14=4; 15=0; 16=4; 17=5; 18=0; 19=0; 20=353; 21=0;
    22=0; 23=0; 24=0; 25=0; 26=5; 27=114; 28=0; 29=1; 30=7; 31=9;
    32=0; 33=2; 34=60; 35=1; 37=1; 38=0; 39=0; 40=6; 41=3306; 42=6;
    43=0; 44=2; 45=0; 46=5; 47=3; 48=1; 49=5; 50=11; 51=6; 52=65;
    53=0; 54=0; 55=28; 56=48; 57=30; 58=0; 59=0; 60=0; 61=1; 62=0;
    63=0; 64=0; 65=0; 66=9; 67=2314; 68=0; 69=1; 70=1; 71=0; 72=389;
    73=5; 74=0; 75=32; 76=0; 77=0; 78=1; 79=1; 80=0; 81=0; 82=0;
    83=1; 84=0; 85=0; 86=2; 87=13; 88=2; 89=25; 90=102; 91=2517; 92=42;
    93=36; 94=0; 95=12; 96=10023; 97=1; 98=52; 99=1; 100=1; 101=0; 102=1;
    103=1; 104=0; 105=0; 106=0; 107=0; 108=8; 109=0; 110=0; 111=0; 112=7;
    113=18; 114=0; 115=27; 116=0; 117=0; 118=0; 119=0; 120=2; 121=16; 122=2;
    123=0; 124=0; 125=1; 126=0; 127=23; 128=0; 129=0; 130=0; 131=0; 132=639;
    133=0; 134=1; 135=1; 136=1; 137=1; 138=4; 139=20; 140=0; 141=11; 142=1;
    143=17; 144=0; 145=3; 146=16; 147=68; 148=18; 149=100; 150=3; 151=1; 152=29;
    153=5; 154=0; 155=0; 156=0; 157=0; 158=22; 159=4; 160=0; 161=10; 162=0;
    163=3; 164=0; 165=0; 166=9; 167=0; 168=2; 169=152; 170=5690; 171=1; 172=1;
    173=0; 174=102; 175=235; 176=0; 177=2; 178=6; 179=0; 180=1; 181=2; 182=6;
    183=108; 184=0; 185=0; 186=8; 187=2; 188=1019; 189=1; 190=0; 191=2191; 192=0;
    193=1; 194=1; 195=0; 196=0; 197=1;
}
foreach ($col in $row91.Keys) {
    $ws.Cells.Item(91, [int]$col).Value = $row91[$col]
}
